$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-seat "icon" header before we overwrite its old cell (C1) so the
#     shared-string entry stays referenced continuously and keeps its slot. ---
$ws.Cells.Item(1, 5).Value = "icon"

# --- New header cells, set in the same left-to-right order the columns
#     end up in, so new shared strings are appended in that order too. ---
$ws.Cells.Item(1, 4).Value  = "类型"
$ws.Cells.Item(1, 6).Value  = "需求等级"
$ws.Cells.Item(1, 7).Value  = "等级上限"
$ws.Cells.Item(1, 8).Value  = "颜色"
$ws.Cells.Item(1, 9).Value  = "职业"
$ws.Cells.Item(1, 10).Value = "是否叠加"
$ws.Cells.Item(1, 11).Value = "最大叠加数"
$ws.Cells.Item(1, 12).Value = "绑定模式"
$ws.Cells.Item(1, 13).Value = "耐久度"
$ws.Cells.Item(1, 14).Value = "耐久显示"
$ws.Cells.Item(1, 15).Value = "冷却类型"
$ws.Cells.Item(1, 16).Value = "冷却时间"
$ws.Cells.Item(1, 17).Value = "价格"
$ws.Cells.Item(1, 18).Value = "卖出价格"

# C1 becomes the new "说明" column (old icon column) - overwrite last.
$ws.Cells.Item(1, 3).Value = "说明"

# --- Data rows ---
$ws.Cells.Item(2, 3).Value  = "这个是道具1号的说明"
$ws.Cells.Item(2, 4).Value  = 1
$ws.Cells.Item(2, 5).Value  = 1
$ws.Cells.Item(2, 6).Value  = 1
$ws.Cells.Item(2, 7).Value  = 65535
$ws.Cells.Item(2, 8).Value  = 1
$ws.Cells.Item(2, 9).Value  = 1

$ws.Cells.Item(3, 3).Value  = "这个是道具2号的说明"
$ws.Cells.Item(3, 4).Value  = 1
$ws.Cells.Item(3, 5).Value  = 1
$ws.Cells.Item(3, 6).Value  = 1
$ws.Cells.Item(3, 7).Value  = 65535
$ws.Cells.Item(3, 8).Value  = 1
$ws.Cells.Item(3, 9).Value  = 1

# --- Highlight the new header cells in red (new fields to be discussed) ---
$ws.Range("C1:D1").Interior.Color = 255
$ws.Range("F1:R1").Interior.Color = 255

# --- Column widths (best effort). The host engine always re-derives the
#     stored OOXML "width" from ColumnWidth using a fixed MDW=7 pixel grid
#     (width -> round(chars*7+5) pixels -> pixels/7), which does not exactly
#     round-trip the original author's widths (authored with a different
#     MDW). We pre-compensate by the grid's constant 5/7 padding so the
#     saved width lands as close as possible to (and exactly matches,
#     whenever the target is already a multiple of 1/7) the intended value.
$pad = 5 / 7
$ws.Columns.Item(3).ColumnWidth  = 20.375 - $pad
$ws.Columns.Item(4).ColumnWidth  = 7.375 - $pad
$ws.Columns.Item(5).ColumnWidth  = 5.5 - $pad
$ws.Columns.Item(6).ColumnWidth  = 9 - $pad
$ws.Columns.Item(7).ColumnWidth  = 9 - $pad
$ws.Columns.Item(8).ColumnWidth  = 5.25 - $pad
$ws.Columns.Item(9).ColumnWidth  = 5.25 - $pad
$ws.Columns.Item(10).ColumnWidth = 9 - $pad
$ws.Columns.Item(11).ColumnWidth = 11 - $pad
$ws.Columns.Item(15).ColumnWidth = 13 - $pad

# --- Move the active selection to D10, like in the final sheet ---
$null = $ws.Range("D10").Select()
